# The source data rows (19-30, excluding the untouched row 25) have had their
# record contents cyclically reshuffled between rows while staying on the
# same row number / position. This reproduces that reshuffle by writing the
# new (post-shuffle) values for every cell that actually changes, including
# re-creating / clearing the "ghost" empty L column cell that travels with
# the "Skogsfru" (Epipogium aphyllum) records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: A, Q, R only
$ws.Range("A19").Value = 111670497
$ws.Range("Q19").Value = 558159.8619213518
$ws.Range("R19").Value = 7068022.886732788

# Row 20: A, Q, R only
$ws.Range("A20").Value = 111671294
$ws.Range("Q20").Value = 558118.4535210516
$ws.Range("R20").Value = 7067742.103054954

# Row 21: becomes the "Skogsfru" / Epipogium aphyllum record (gains L21)
$ws.Range("A21").Value = 111670477
$ws.Range("B21").Value = 96346
$ws.Range("E21").Value = 620
$ws.Range("F21").Value = "Skogsfru"
$ws.Range("G21").Value = "Epipogium aphyllum"
$ws.Range("H21").Value = "Sw."
$ws.Range("L21").Font.Bold = $false
$ws.Range("Q21").Value = 558155.0815836267
$ws.Range("R21").Value = 7068017.481975557

# Row 22: A, Q, R only
$ws.Range("A22").Value = 111670567
$ws.Range("Q22").Value = 558129.9933989302
$ws.Range("R22").Value = 7067958.536170656

# Row 23: becomes "Stuplav" / Nephroma bellum
$ws.Range("A23").Value = 111671188
$ws.Range("B23").Value = 78605
$ws.Range("E23").Value = 6462
$ws.Range("F23").Value = "Stuplav"
$ws.Range("G23").Value = "Nephroma bellum"
$ws.Range("H23").Value = "(Spreng.) Tuck."

# Row 24: becomes "Skrovellav" / Lobaria scrobiculata (loses L24)
$ws.Range("A24").Value = 111671201
$ws.Range("B24").Value = 78579
$ws.Range("E24").Value = 2081
$ws.Range("F24").Value = "Skrovellav"
$ws.Range("G24").Value = "Lobaria scrobiculata"
$ws.Range("H24").Value = "(Scop.) DC."
$ws.Range("L24").ClearContents()
$ws.Range("Q24").Value = 558250.1783714101
$ws.Range("R24").Value = 7067936.828089682

# Row 26: becomes "Bårdlav" / Nephroma parile (loses L26)
$ws.Range("A26").Value = 111671190
$ws.Range("B26").Value = 78611
$ws.Range("D26").Value = "LC"
$ws.Range("E26").Value = 6463
$ws.Range("F26").Value = "Bårdlav"
$ws.Range("G26").Value = "Nephroma parile"
$ws.Range("H26").Value = "(Ach.) Ach."
$ws.Range("L26").ClearContents()
$ws.Range("Q26").Value = 558215.9329796816
$ws.Range("R26").Value = 7067869.292590594

# Row 27: A, Q, R only
$ws.Range("A27").Value = 111671226
$ws.Range("Q27").Value = 558118.4535210516
$ws.Range("R27").Value = 7067742.103054954

# Row 28: becomes the "Skogsfru" / Epipogium aphyllum record (gains L28)
$ws.Range("A28").Value = 111670510
$ws.Range("B28").Value = 96346
$ws.Range("E28").Value = 620
$ws.Range("F28").Value = "Skogsfru"
$ws.Range("G28").Value = "Epipogium aphyllum"
$ws.Range("H28").Value = "Sw."
$ws.Range("L28").Font.Bold = $false
$ws.Range("Q28").Value = 558124.4538526792
$ws.Range("R28").Value = 7067994.321708324

# Row 29: becomes "Lunglav" / Lobaria pulmonaria (loses L29)
$ws.Range("A29").Value = 111671179
$ws.Range("B29").Value = 78578
$ws.Range("E29").Value = 6458
$ws.Range("F29").Value = "Lunglav"
$ws.Range("G29").Value = "Lobaria pulmonaria"
$ws.Range("H29").Value = "(L.) Hoffm."
$ws.Range("L29").ClearContents()
$ws.Range("Q29").Value = 558215.9656782644
$ws.Range("R29").Value = 7067867.520903144

# Row 30: becomes the "Skogsfru" / Epipogium aphyllum record (gains L30)
$ws.Range("A30").Value = 111670558
$ws.Range("B30").Value = 96346
$ws.Range("D30").Value = "NT"
$ws.Range("E30").Value = 620
$ws.Range("F30").Value = "Skogsfru"
$ws.Range("G30").Value = "Epipogium aphyllum"
$ws.Range("H30").Value = "Sw."
$ws.Range("L30").Font.Bold = $false
$ws.Range("Q30").Value = 558133.6011735104
$ws.Range("R30").Value = 7067979.426396712
